$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B4 - same content/position, wording unchanged (kept for completeness/clarity)
$ws.Range("B4").Value = "Add a description for each prey code if desired (taxa name or general descriptor)"

# Update B8 - sequential integer instructions now describe FIRST row = UNID (TypeN = 0)
$ws.Range("B8").Value = "Add a sequential integer for each unique prey type. Note that the FIRST row is reserved for the `"UNID`" prey type, with TypeN value of 0. Begin adding prey types on subsequent rows, with TypeN values of 1, 2… etc."

# Update B9 - prey type list instructions now describe FIRST row = UNID, and 3-4 letter code
$ws.Range("B9").Value = "Enter a list of mutually exclusive prey types for analysis, with a unique 3 or 4 letter code for each. You should limit the number of possible prey types such that each type is represented by a reasonable number of instances in the data set (see `"N`" nolumn in previous worksheet: a minimum of 25 is a good rule of thumb, at least 50 preferred).  Thus prey codes with few observations should be combined into a single prey type.  NOTE: The FIRST row is reserved for the `"UNID`" prey type (un-identified prey): leave this as-is."

# Update B10 and B11 - content unchanged, kept for completeness/consistency
$ws.Range("B10").Value = "Select a prey class for each prey type: refer to list in next worksheet for a complete list of possible prey classes. Leave this blank for `"UNID`" prey type"
$ws.Range("B11").Value = "These are the codes for the available species from prey database having both size-mass and caloric data"

# Move the active cell selection from B5 to B2
$ws.Range("B2").Select()
